# Auto-generated edit script: updates computed market-price / profit
# columns (H-N) across all 8 job sheets per the scheduled-runner diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 1073.6666
$ws.Range("I8").Value = 1073.6666
$ws.Range("K8").Value = 3220.9998
$ws.Range("M8").Value = -3081.9998
$ws.Range("H64").Value = 6154.8184
$ws.Range("I64").Value = 6212.5
$ws.Range("K64").Value = 6212.5
$ws.Range("M64").Value = -5964.5
$ws.Range("H67").Value = 6154.8184
$ws.Range("I67").Value = 6212.5
$ws.Range("K67").Value = 6212.5
$ws.Range("M67").Value = -5354.5
$ws.Range("H111").Value = 2943085.2
$ws.Range("I111").Value = 3012.4285
$ws.Range("K111").Value = 9037.2855
$ws.Range("M111").Value = -5970.2855

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2244.7778
$ws.Range("I45").Value = 2073.7334
$ws.Range("J45").Value = 3100
$ws.Range("K45").Value = 2073.7334
$ws.Range("L45").Value = 3100
$ws.Range("M45").Value = -1696.7334
$ws.Range("N45").Value = -3854
$ws.Range("H61").Value = 4251.387
$ws.Range("I61").Value = 3260.2856
$ws.Range("K61").Value = 3260.2856
$ws.Range("M61").Value = -3048.2856
$ws.Range("H132").Value = 6859.3335
$ws.Range("I132").Value = 4931.2
$ws.Range("J132").Value = 16500
$ws.Range("K132").Value = 14793.6
$ws.Range("L132").Value = 49500
$ws.Range("M132").Value = -12263.6
$ws.Range("N132").Value = -54560
$ws.Range("H136").Value = 4251.387
$ws.Range("I136").Value = 3260.2856
$ws.Range("K136").Value = 9780.856800000001
$ws.Range("M136").Value = -7230.856800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4821.1113
$ws.Range("I20").Value = 4455.75
$ws.Range("J20").Value = 5113.4
$ws.Range("K20").Value = 4455.75
$ws.Range("L20").Value = 5113.4
$ws.Range("M20").Value = -4208.75
$ws.Range("N20").Value = -5607.4
$ws.Range("H107").Value = 2101.8928
$ws.Range("I107").Value = 1942.8334
$ws.Range("J107").Value = 3056.25
$ws.Range("K107").Value = 1942.8334
$ws.Range("L107").Value = 3056.25
$ws.Range("M107").Value = -22.83339999999998
$ws.Range("N107").Value = -6896.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5161.3076
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 5161.3076
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 5161.3076
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -5735.3076
$ws.Range("H99").Value = 2662.125
$ws.Range("I99").Value = 2659.4
$ws.Range("J99").Value = 2666.6667
$ws.Range("K99").Value = 2659.4
$ws.Range("L99").Value = 2666.6667
$ws.Range("M99").Value = -1161.4
$ws.Range("N99").Value = -5662.6667
$ws.Range("H105").Value = 1072.7142
$ws.Range("I105").Value = 1072.7142
$ws.Range("K105").Value = 1072.7142
$ws.Range("M105").Value = 674.2858000000001
$ws.Range("H113").Value = 5161.3076
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 5161.3076
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 5161.3076
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -9501.3076
$ws.Range("H126").Value = 2662.125
$ws.Range("I126").Value = 2659.4
$ws.Range("J126").Value = 2666.6667
$ws.Range("K126").Value = 7978.200000000001
$ws.Range("L126").Value = 8000.000100000001
$ws.Range("M126").Value = -5508.200000000001
$ws.Range("N126").Value = -12940.0001
$ws.Range("H135").Value = 120000
$ws.Range("J135").Value = 120000
$ws.Range("L135").Value = 120000
$ws.Range("N135").Value = -130140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 85
$ws.Range("I38").Value = 85
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 255
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 92
$ws.Range("N38").ClearContents()
$ws.Range("H64").Value = 7775.846
$ws.Range("I64").Value = 9996.333000000001
$ws.Range("J64").Value = 7109.7
$ws.Range("K64").Value = 29988.999
$ws.Range("L64").Value = 21329.1
$ws.Range("M64").Value = -29718.999
$ws.Range("N64").Value = -21869.1
$ws.Range("H67").Value = 7775.846
$ws.Range("I67").Value = 9996.333000000001
$ws.Range("J67").Value = 7109.7
$ws.Range("K67").Value = 29988.999
$ws.Range("L67").Value = 21329.1
$ws.Range("M67").Value = -29052.999
$ws.Range("N67").Value = -23201.1
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H113").Value = 98
$ws.Range("J113").Value = 98
$ws.Range("L113").Value = 294
$ws.Range("N113").Value = -4634
$ws.Range("H117").Value = 1890
$ws.Range("J117").Value = 1890
$ws.Range("L117").Value = 5670
$ws.Range("N117").Value = -12554
$ws.Range("H133").Value = 0
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6942.148
$ws.Range("I122").Value = 7673.05
$ws.Range("J122").Value = 4853.857
$ws.Range("K122").Value = 23019.15
$ws.Range("L122").Value = 14561.571
$ws.Range("M122").Value = -20569.15
$ws.Range("N122").Value = -19461.571
$ws.Range("H132").Value = 4794.207
$ws.Range("I132").Value = 2515.8235
$ws.Range("K132").Value = 7547.470499999999
$ws.Range("M132").Value = -5017.470499999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2187.5
$ws.Range("I16").Value = 1306.6364
$ws.Range("K16").Value = 1306.6364
$ws.Range("M16").Value = -1136.6364
$ws.Range("H46").Value = 6075
$ws.Range("I46").Value = 1858.3334
$ws.Range("J46").Value = 7340
$ws.Range("K46").Value = 1858.3334
$ws.Range("L46").Value = 7340
$ws.Range("M46").Value = -1670.3334
$ws.Range("H115").Value = 106498.5
$ws.Range("J115").Value = 106498.5
$ws.Range("L115").Value = 106498.5
$ws.Range("N115").Value = -108848.5
$ws.Range("H125").Value = 181249.19
$ws.Range("J125").Value = 181249.19
$ws.Range("L125").Value = 181249.19
$ws.Range("N125").Value = -191089.19
$ws.Range("H132").Value = 3718.6924
$ws.Range("I132").Value = 3032.15
$ws.Range("K132").Value = 9096.450000000001
$ws.Range("M132").Value = -6566.450000000001
$ws.Range("H133").Value = 50305.25
$ws.Range("J133").Value = 50305.25
$ws.Range("L133").Value = 50305.25
$ws.Range("N133").Value = -55365.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 32665.666
$ws.Range("J34").Value = 31999.5
$ws.Range("L34").Value = 31999.5
$ws.Range("N34").Value = -32405.5
$ws.Range("H43").Value = 120000
$ws.Range("J43").Value = 120000
$ws.Range("L43").Value = 120000
$ws.Range("N43").Value = -120298
$ws.Range("H103").Value = 71293
$ws.Range("J103").Value = 71293
$ws.Range("L103").Value = 71293
$ws.Range("N103").Value = -73637
$ws.Range("H126").Value = 3979.6128
$ws.Range("J126").Value = 8748.375
$ws.Range("L126").Value = 26245.125
$ws.Range("N126").Value = -31185.125
